# Changes for Service - 13 March 2023
# Applies the additions to the "Constants" sheet: new ServiceReportFolderPath,
# ServiceReportFileExtension and LocalInvoiceFolder config rows, and renames
# the "Services" value from "Services" to "Service".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# Step 1: insert a new row above row 36 ("RootFolder" and everything below
# shifts down by one) for the new ServiceReportFolderPath setting.
# ---------------------------------------------------------------------------
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).ClearFormats()
$ws.Rows.Item(36).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Step 2: after the previous insert, "TempTracker" now sits on row 42.
# Insert three new rows right after it (at what is currently row 43, the
# blank spacer row) for LocalInvoiceFolder, a blank spacer, and
# ServiceReportFileExtension.
# ---------------------------------------------------------------------------
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).Insert()

$ws.Rows.Item(43).ClearFormats()
$ws.Rows.Item(43).RowHeight = 14.25

$ws.Rows.Item(44).ClearFormats()
$ws.Rows.Item(44).RowHeight = 14.25

$ws.Rows.Item(45).ClearFormats()
$ws.Rows.Item(45).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Now populate the cell values. The order in which distinct new text values
# are written determines the order new shared strings are appended in, so
# we follow the same order as the source workbook's revision.
# ---------------------------------------------------------------------------

# Row 36: ServiceReportFolderPath
$ws.Cells.Item(36, 2).Value = "C:\Box\Personal_603214\External\RAJAT\Development\GF Service Report"
$ws.Cells.Item(36, 1).Value = "ServiceReportFolderPath"

# Row 45: ServiceReportFileExtension
$ws.Cells.Item(45, 2).Value = ".pdf"
$ws.Cells.Item(45, 1).Value = "ServiceReportFileExtension"

# Row 49 (the "Services" row): rename the value from "Services" to "Service"
$ws.Cells.Item(49, 2).Value = "Service"

# Row 43: LocalInvoiceFolder
$ws.Cells.Item(43, 1).Value = "LocalInvoiceFolder"
$ws.Cells.Item(43, 2).Value = "Data\Output\Invoices"
$ws.Cells.Item(43, 3).Value = "Output"

# ---------------------------------------------------------------------------
# Step 4: the GFHomePage_URL row (now row 52) loses its wrap text and the
# row height shrinks from 58 to 14.5.
# ---------------------------------------------------------------------------
$ws.Cells.Item(52, 2).WrapText = $false
$ws.Rows.Item(52).RowHeight = 14.5

Write-Host "Applied Service configuration changes to Constants sheet."
